# "corrected times to UTC" — shift every match kickoff time in the
# "Matches" sheet (column D, rows 2-65) back by 10 hours (10/24 of a day),
# converting the stored local kickoff times to UTC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matches")

$hoursOffset = 10 / 24

for ($r = 2; $r -le 65; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value2 = $cell.Value2 - $hoursOffset
}

# Reflect where the author's selection ended up after the edit: sheet
# scrolled back to the top with D2 selected (rather than parked at B65
# after editing the last row).
$ws.Activate()
$ws.Range("D2").Select()
